$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------------
# "Weekly Progress" sheet: a new week's numbers were recorded, so a new
# top data row is inserted above the existing history (everything else
# shifts down by one row).
# --------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Weekly Progress")

# Insert a new row at row 2 - this pushes all the existing weekly rows
# down by one and keeps the column C "=SUM(Bn,-Dn)" formulas correct.
$ws5.Rows.Item(2).Insert()

# Bring over the date-column number formatting from the row that used to
# be on top (now row 3) so the new row looks like the rest of the table.
$ws5.Range("A3:D3").Copy()
$ws5.Range("A2:D2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in this week's numbers.
$ws5.Range("A2").Value = 43735
$ws5.Range("B2").Value = 273
$ws5.Range("D2").Value = 44
$ws5.Range("C2").Formula = "=SUM(B2,-D2)"

# Grow the "Table15" list object so it covers the new row too.
$lo = $ws5.ListObjects.Item(1)
$lo.Resize($ws5.Range("A1:D19"))

# Leave the selection where the author left it while adding the new row.
$ws5.Range("A3").Select()

# --------------------------------------------------------------------------
# "Areas Features Validations" sheet: the author had scrolled further down
# the sheet (the active tab / selection moved accordingly).
# --------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Areas Features Validations")
$ws2.Activate()
$ws2.Range("A159").Select()
$excel.ActiveWindow.ScrollRow = 159
$ws2.Range("D41").Select()
